# Insert a new daily record at row 15 (Fruta / hortaliza, semanal update).
# Excel shifts the existing rows 15..51 down to 16..52, preserving their
# data and formatting (the date column's numeric format carries with it).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(15).Insert()

# Populate the newly-inserted (and now blank) row 15 with the new record.
$ws.Cells.Item(15, 1).Value  = 1
$ws.Cells.Item(15, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(15, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(15, 4).Value  = 44600
$ws.Cells.Item(15, 5).Value  = 15
$ws.Cells.Item(15, 6).Value  = "Fruta"
$ws.Cells.Item(15, 7).Value  = 100103
$ws.Cells.Item(15, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(15, 9).Value  = 100103006
$ws.Cells.Item(15, 10).Value = "Nectarín"
$ws.Cells.Item(15, 11).Value = "June Pearl"
$ws.Cells.Item(15, 12).Value = "Segunda"
$ws.Cells.Item(15, 13).Value = 250
$ws.Cells.Item(15, 14).Value = 17000
$ws.Cells.Item(15, 15).Value = 18000
$ws.Cells.Item(15, 16).Value = 17500
$ws.Cells.Item(15, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(15, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(15, 19).Value = 972
$ws.Cells.Item(15, 20).Value = 18
